# Fixed README.md stats and docx preparation for all Renaissance - JDK 17 -
# Shenandoah GC tests.
#
# The single-column results table had its "per-iteration" rows (6 through
# ~12) holding stale placeholder values, while the true per-iteration
# numbers were packed (tab separated) into the summary-row runs near the
# bottom of the table. This moves those numbers into their own rows and
# collapses the summary rows back down to a single value each.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Simple single-run value corrections (rows 1-6) ---
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "1098"
$t.Cell(5, 1).Range.Text = "0.00002"
$t.Cell(6, 1).Range.Text = "0.00290"

# --- Remove the three stale placeholder rows that followed (originally
#     rows 7, 8, 9 holding 0.00013 / 0.00004 / 0.00015). Deleting row 7
#     three times removes each in turn as the rest shift up. ---
$t.Rows(7).Delete()
$t.Rows(7).Delete()
$t.Rows(7).Delete()

# Row 7 (0.00020) is unchanged; rows 8 and 9 get corrected values.
$t.Cell(8, 1).Range.Text = "0.00010"
$t.Cell(9, 1).Range.Text = "0.00028"

# --- Insert three new rows (for 0.00038 / 0.00069 / 0.26546) right before
#     the "100.0" row that currently sits at index 10. ---
$beforeRow = $t.Rows(10)
$t.Rows.Add($beforeRow)
$t.Rows.Add($beforeRow)
$t.Rows.Add($beforeRow)
$t.Cell(10, 1).Range.Text = "0.00038"
$t.Cell(11, 1).Range.Text = "0.00069"
$t.Cell(12, 1).Range.Text = "0.26546"

# --- Collapse the three tab-separated summary rows at the bottom of the
#     table down to their single leading value now that the rest of the
#     figures live in their own rows above. ---
$t.Cell(44, 1).Range.Text = "99.81"
$t.Cell(45, 1).Range.Text = "0.27"
$t.Cell(46, 1).Range.Text = "139"
